$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." built from runs:
#   "Versi" | "on" | " 2" | "." (with a spellcheck wrapper around "Versi"+"on"
#   and a _GoBack bookmark sitting between " 2" and ".").
# Target reading: "Version 1." built from runs "Version" | " 1." (the
# trailing "." run disappears, merged into the " 1." run).

# 1) Drop the trailing "." run (rightmost edit first so earlier offsets
#    used below stay valid).
$d.Range(9, 10).Text = ""

# 2) Turn " 2" into " 1.".
$d.Range(7, 9).Text = " 1."

# 3) Merge the "Versi" + "on" runs into a single "Version" run. The text
#    is already "Version", so a same-text assignment would be a no-op;
#    nudge it through a temporary value first to force the run merge.
$d.Range(0, 7).Text = "Version_"
$d.Range(0, 8).Text = "Version"
